$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the booking JSON sample text: "date" -> "dateRequired" ---
# The placeholder JSON is reused (identical text) by D27, E26, E27, D29 and E29
# (all "booking" sample params / sample return cells), so every cell that
# currently shows it needs to be updated to keep them in sync.
$lcurly = [string][char]0x201C   # “
$rcurly = [string][char]0x201D   # ”
$quote  = [string][char]0x22     # "

$part1 = "{" + $rcurly + "roomID" + $quote + ":" + $rcurly + "xxx" + $rcurly + ","
$part2 = $lcurly + "dateRequired" + $rcurly + ":" + $rcurly + "xxx" + $rcurly + ","
$part3 = $lcurly + "userName:" + $rcurly + "xxx" + $rcurly + ","
$part4 = $rcurly + "reason" + $rcurly + ":" + $rcurly + "xxxx" + $rcurly + "  }"
$bookingJson = $part1 + $part2 + $part3 + $part4

$ws.Range("E26").Value = $bookingJson
$ws.Range("D27").Value = $bookingJson
$ws.Range("E27").Value = $bookingJson
$ws.Range("D29").Value = $bookingJson
$ws.Range("E29").Value = $bookingJson

# --- 2. E26 / E27 switch from general to justified alignment (matches the
#        style already used by the rest of the table, making the old
#        general-alignment style redundant) ---
$ws.Range("E26").HorizontalAlignment = -4130
$ws.Range("E27").HorizontalAlignment = -4130

# --- 3. Row height tweaks in the "Proposed Solutions" table ---
$ws.Rows(20).RowHeight = 13.25
$ws.Rows(21).RowHeight = 13.25
$ws.Rows(22).RowHeight = 13.25
$ws.Rows(23).RowHeight = 13.25
$ws.Rows(24).RowHeight = 13.25
$ws.Rows(25).RowHeight = 13.25
$ws.Rows(26).RowHeight = 25.3
$ws.Rows(27).RowHeight = 25.3
$ws.Rows(28).RowHeight = 13.25

# --- 4. Move the active selection to A24 ---
$ws.Range("A24").Select() | Out-Null
